# The workbook documents a DB-normalization exercise. On Sheet4 ("영화관
# 예매 시스템" / movie-theater reservation ERD) the author reworked the
# schema: the old "가격(prices) / Price_by_seat / movie_screening" tables
# were replaced with "show_time / reserve_seat" tables. Update the cells
# that describe those entities/attributes accordingly.

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)

# Cells that keep using already-existing text (just re-typed / untouched):
$ws4.Range("E24").Value = "show_time"
$ws4.Range("F25").Value = "show_time"
$ws4.Range("G25").Value = "movie_id"
$ws4.Range("H25").Value = "screen_id"
$ws4.Range("L25").Value = "prices"

# M25 is a brand-new cell; give it the same formatting as its neighbour
# L25 before filling in its value.
$ws4.Range("L25").Copy()
$ws4.Range("M25").PasteSpecial(-4122)  # xlPasteFormats
$ws4.Range("M25").Value = "reserve_id"
$excel.CutCopyMode = $false

# Cells introducing genuinely new vocabulary - entered in the same order the
# author typed them in (matches the order new entries were appended to the
# shared-string table):
$ws4.Range("E25").Value = "show_time_id"
$ws4.Range("F24").Value = "상영일정"
$ws4.Range("G29").Value = "예매"
$ws4.Range("J24").Value = "reserve_seat"
$ws4.Range("J25").Value = "reserve_seat_id"
$ws4.Range("K25").Value = "seat"
$ws4.Range("K24").Value = "예매 상세"

# movie_screening_id -> show_time_id (reuses the string created above)
$ws4.Range("H30").Value = "show_time_id"

# The old "prices"/"가격"/"Price_by_seat_id"/"prices_id" reference cells no
# longer apply to the reworked schema - remove them outright.
$ws4.Range("L29").Clear()
$ws4.Range("M29").Clear()
$ws4.Range("I30").Clear()
$ws4.Range("L30").Clear()
$ws4.Range("M30").Clear()

# --- Restore the cursor position recorded in the saved view ----------------
[void]$ws4.Range("H23").Select()
